# chore: adapt column header formatting to respective input file names
#
# Renames the header row of the "AHB-Diff" sheet so the "_old"/"_new"
# suffixes become "_FV2410"/"_FV2504" (the actual format-version names of
# the two compared input files), wraps the data range A1:U62 in an Excel
# Table ("Table1") using the same (new) header names as its column names,
# and freezes the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------
$oldSuffixHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$col = 1
foreach ($base in $oldSuffixHeaders) {
    $ws.Cells.Item(1, $col).Value = "$base`_FV2410"
    $col++
}

# column K ("diff") is unchanged
$col++

foreach ($base in $oldSuffixHeaders) {
    $ws.Cells.Item(1, $col).Value = "$base`_FV2504"
    $col++
}

# --- 2. Turn the range into a proper Excel Table -------------------------
$tableRange = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
